# Add the "alu" result block (rows 14-16) to the results sheet and move
# the active selection, per commit "add alu result and change lib for
# xor from 1000 to 5".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is a blank spacer row that only carries the bold/right-hand
# column formatting in U (matches the style used for the other spacer
# rows, e.g. row 7, 10, 11).
$ws.Range("U7").Copy()
$ws.Range("U14").PasteSpecial(-4122)

# Row 15: new "alu" size=4 entry.
$ws.Range("A15").Value = "alu"
$ws.Range("B15").Value = 4
$ws.Range("P15").Value = 26
$ws.Range("Q15").Value = 81
$ws.Range("R15").Formula = "=SUM(P15:Q15)"
$ws.Range("S15").Value = 163
$ws.Range("T15").Formula = "=SUM(P15:R15)"
$ws.Range("U15").Formula = "=SUM(P15:Q15)+5*S15"

# Row 16: new "alu" size=8 entry (continuation row, no label in A).
$ws.Range("B16").Value = 8
$ws.Range("P16").Value = 139
$ws.Range("Q16").Value = 220
$ws.Range("R16").Formula = "=SUM(P16:Q16)"
$ws.Range("S16").Value = 422
$ws.Range("T16").Formula = "=SUM(P16:R16)"
$ws.Range("U16").Formula = "=SUM(P16:Q16)+5*S16"

# Match the formatting used by the other "Function" label / total cells.
$ws.Range("A12").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("U13").Copy()
$ws.Range("U15:U16").PasteSpecial(-4122)

# Move the active cell/selection as recorded in the saved view state.
$ws.Range("U20").Select()
